{"js": "// The commit moves the contract's dates \"1 day earlier\" relative to\n// today's mailing (09 Sep -> 10 Sep) and fixes the term-start date to\n// be \"after the due date\" (01 Oct -> 11 Sep).\n//\n// There are 3 literal occurrences of \"09 September 2025\" in the body\n// (the offer-date sentence, the signature-block date, and the\n// \"CONFIRM\" date line) that all become \"10 September 2025\", plus one\n// occurrence of \"01 October 2025\" (the agreement term start) that\n// becomes \"11 September 2025\".\n\nconst body = context.document.body;\n\n// 1) \"09 September 2025\" -> \"10 September 2025\" (all matches).\nconst offerDateHits = body.search(\"09 September 2025\", { matchCase: true });\nofferDateHits.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < offerDateHits.items.length; i++) {\n  offerDateHits.items[i].insertText(\"10 September 2025\", \"Replace\");\n}\nawait context.sync();\n\n// 2) \"01 October 2025\" -> \"11 September 2025\" (term commencement date).\nconst termDateHits = body.search(\"01 October 2025\", { matchCase: true });\ntermDateHits.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < termDateHits.items.length; i++) {\n  termDateHits.items[i].insertText(\"11 September 2025\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# The commit moves the contract's dates \"1 day earlier\" relative to\n# today's mailing (09 Sep -> 10 Sep) and fixes the term-start date to\n# be \"after the due date\" (01 Oct -> 11 Sep).\n#\n# There are 3 literal occurrences of \"09 September 2025\" in the body\n# (the offer-date sentence, the signature-block date, and the\n# \"CONFIRM\" date line) that all become \"10 September 2025\", plus one\n# occurrence of \"01 October 2025\" (the agreement term start) that\n# becomes \"11 September 2025\".\n\n$d = $word.ActiveDocument\n\n# 1) \"09 September 2025\" -> \"10 September 2025\" (all matches).\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"09 September 2025\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"10 September 2025\"\n$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 2) \"01 October 2025\" -> \"11 September 2025\" (term commencement date).\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"01 October 2025\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"11 September 2025\"\n$find2.Execute($find2.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n"}
